$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 75
$ws.Range("A75").Value = "05/01/2026 11:58:11"
$ws.Range("B75").Value = "05/01 11:54"
$ws.Range("C75").Value = "Folha de S.Paulo - Mercado - Principal"
$ws.Range("D75").Value = "Pagamento do IPVA 2026 já está liberado em bancos; veja onde quitar e calendário de São Paulo"
$ws.Range("E75").Value = "https://redir.folha.com.br/redir/online/mercado/rss091/*https://www1.folha.uol.com.br/mercado/2026/01/pagamento-do-ipva-2026-ja-esta-liberado-em-bancos-veja-onde-quitar-e-calendario-de-sao-paulo.shtml"
$ws.Range("F75").Value = "imposto"
$ws.Range("G75").Value = " da Fazenda e Planejamento do Estado de São Paulo) informou que liberou as informações do &lt;b&gt;imposto&lt;/b&gt; às instituições financeiras, e que a oferta do serviço aos clientes depende de cada banco"

# Row 76
$ws.Range("A76").Value = "05/01/2026 11:58:12"
$ws.Range("B76").Value = "05/01 11:40"
$ws.Range("C76").Value = "Folha de S.Paulo - Mercado - Principal"
$ws.Range("D76").Value = "UE diz que houve avanços para fechar acordo 'em breve' com Mercosul"
$ws.Range("E76").Value = "https://redir.folha.com.br/redir/online/mercado/rss091/*https://www1.folha.uol.com.br/mercado/2026/01/ue-diz-que-houve-avancos-para-fechar-acordo-em-breve-com-mercosul.shtml"
$ws.Range("F76").Value = "comissão"
$ws.Range("G76").Value = "A Comissão Europeia afirmou nesta segunda-feira (5) que houve ""avanços"" entre os Estados europeus em"

# Row 77
$ws.Range("A77").Value = "05/01/2026 11:58:13"
$ws.Range("B77").Value = "05/01 11:34"
$ws.Range("C77").Value = "Folha de S.Paulo - Mercado - Principal"
$ws.Range("D77").Value = "Nova tabela do Imposto de Renda começa a valer; veja o que muda nos salários"
$ws.Range("E77").Value = "https://redir.folha.com.br/redir/online/mercado/rss091/*https://www1.folha.uol.com.br/mercado/2026/01/nova-tabela-do-imposto-de-renda-comeca-a-valer-veja-o-que-muda-nos-salarios.shtml"
$ws.Range("F77").Value = "imposto"
$ws.Range("G77").Value = "A nova tabela do &lt;a href=""https://www1.folha.uol.com.br/folha-topicos/&lt;b&gt;imposto&lt;/b&gt;-de-renda/""&gt;Imposto de Renda&lt;/a&gt; começa a valer em 1º de janeiro de 2026 e zera a cobrança"

# Row 78
$ws.Range("A78").Value = "05/01/2026 11:58:14"
$ws.Range("B78").Value = "05/01 09:36"
$ws.Range("C78").Value = "g1 > Economia"
$ws.Range("D78").Value = "França vai suspender importação de frutas do Mercosul com agrotóxicos proibidos na Europa"
$ws.Range("E78").Value = "https://g1.globo.com/economia/agronegocios/noticia/2026/01/05/franca-suspende-importacao-de-frutas-do-mercosul-com-agrotoxicos-proibidos-na-europa.ghtml"
$ws.Range("F78").Value = "comissão"
$ws.Range("G78").Value = " entre os dois blocos (UE e Mercosul). Ele havia sido fechado em dezembro de 2024 entre a Comissão Europeia, o órgão executivo da UE, com Argentina, Brasil, Paraguai e Uruguai. `nAgricultor"

